$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 43.63
$ws.Range("C3").Value = 28.71
$ws.Range("C4").Value = 5.71
$ws.Range("C5").Value = 7.14
